$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Smoke"
$ws.Name = "Smoke"

# Update row 2 values: E2/F2/G2 now hold different product/test data
$ws.Range("E2").Value = "Paracip"
$ws.Range("F2").Value = "B"
$ws.Range("G2").Value = "C:\Users\DELL\Desktop\\download (2).jpg"

# Widen column E to a fixed custom width of 25 (was auto bestFit ~18.29)
$ws.Columns("E").ColumnWidth = 24.166666666666668

# Move the active selection to G13
$ws.Range("G13").Select() | Out-Null
